# A new weekly price record needs to be inserted at row 17 (pushing the
# existing rows 17-49 down to 18-50, so the last existing row ends up as
# the new row 50). Row 17 then gets populated with its own fresh values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17; Excel shifts rows
# 17:49 down to 18:50 and extends the used range automatically.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new record.
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44481
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112038
$ws.Range("G17").Value = "Cebollín baby"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = 950
$ws.Range("N17").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 475
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = "Hortaliza"
